$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.174496644295302
$ws.Range("C2").Value = 0.5939597315436241
$ws.Range("J2").Value = 0.003355704697986577
$ws.Range("P2").Value = 0.1174496644295302
$ws.Range("S2").Value = 0.1107382550335571
$ws.Range("C3").Value = 0.005524861878453038
$ws.Range("J3").Value = 0.01657458563535912
$ws.Range("P3").Value = 0.7182320441988951
$ws.Range("S3").Value = 0.2596685082872928
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.6458333333333334
$ws.Range("S4").Value = 0.2916666666666667
$ws.Range("B6").Value = 0.06
$ws.Range("D6").Value = 0.025
$ws.Range("F6").Value = 0.07000000000000001
$ws.Range("J6").Value = 0.24
$ws.Range("O6").Value = 0.015
$ws.Range("Q6").Value = 0.185
$ws.Range("R6").Value = 0.055
$ws.Range("S6").Value = 0.35
$ws.Range("B7").Value = 0.08641975308641975
$ws.Range("D7").Value = 0.01851851851851852
$ws.Range("F7").Value = 0.03703703703703703
$ws.Range("J7").Value = 0.1419753086419753
$ws.Range("O7").Value = 0.02469135802469136
$ws.Range("Q7").Value = 0.2407407407407407
$ws.Range("R7").Value = 0.05555555555555555
$ws.Range("S7").Value = 0.3950617283950617
$ws.Range("B8").Value = 0.1167883211678832
$ws.Range("D8").Value = 0.0194647201946472
$ws.Range("E8").Value = 0.0024330900243309
$ws.Range("F8").Value = 0.04866180048661801
$ws.Range("J8").Value = 0.1265206812652068
$ws.Range("O8").Value = 0.0218978102189781
$ws.Range("Q8").Value = 0.1897810218978102
$ws.Range("R8").Value = 0.06082725060827251
$ws.Range("S8").Value = 0.413625304136253
$ws.Range("B9").Value = 0.08333333333333333
$ws.Range("D9").Value = 0.02976190476190476
$ws.Range("F9").Value = 0.07738095238095238
$ws.Range("J9").Value = 0.1726190476190476
$ws.Range("O9").Value = 0.01785714285714286
$ws.Range("Q9").Value = 0.2023809523809524
$ws.Range("R9").Value = 0.05357142857142857
$ws.Range("S9").Value = 0.3630952380952381
$ws.Range("B10").Value = 0.116030534351145
$ws.Range("D10").Value = 0.02213740458015267
$ws.Range("E10").Value = 0.001526717557251908
$ws.Range("F10").Value = 0.06946564885496183
$ws.Range("J10").Value = 0.133587786259542
$ws.Range("O10").Value = 0.01755725190839695
$ws.Range("Q10").Value = 0.2297709923664122
$ws.Range("R10").Value = 0.06030534351145038
$ws.Range("S10").Value = 0.349618320610687
$ws.Range("F11").Value = 0.003460207612456748
$ws.Range("G11").Value = 0.1384083044982699
$ws.Range("J11").Value = 0.09688581314878893
$ws.Range("K11").Value = 0.2041522491349481
$ws.Range("L11").Value = 0.5467128027681661
$ws.Range("S11").Value = 0.01038062283737024
$ws.Range("G12").Value = 0.6848484848484848
$ws.Range("J12").Value = 0.2242424242424242
$ws.Range("K12").Value = 0.01212121212121212
$ws.Range("L12").Value = 0.03636363636363636
$ws.Range("S12").Value = 0.04242424242424243
$ws.Range("G13").Value = 0.64
$ws.Range("J13").Value = 0.32
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.02369668246445497
$ws.Range("H15").Value = 0.1137440758293839
$ws.Range("I15").Value = 0.07582938388625593
$ws.Range("J15").Value = 0.4170616113744076
$ws.Range("K15").Value = 0.07582938388625593
$ws.Range("M15").Value = 0.01421800947867299
$ws.Range("O15").Value = 0.03317535545023697
$ws.Range("S15").Value = 0.2464454976303317
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.1648936170212766
$ws.Range("I16").Value = 0.101063829787234
$ws.Range("J16").Value = 0.4680851063829787
$ws.Range("K16").Value = 0.09574468085106383
$ws.Range("M16").Value = 0.005319148936170213
$ws.Range("O16").Value = 0.04787234042553191
$ws.Range("S16").Value = 0.09574468085106383
$ws.Range("F17").Value = 0.004106776180698152
$ws.Range("H17").Value = 0.2032854209445585
$ws.Range("I17").Value = 0.06570841889117043
$ws.Range("J17").Value = 0.4373716632443532
$ws.Range("K17").Value = 0.08624229979466119
$ws.Range("M17").Value = 0.01026694045174538
$ws.Range("O17").Value = 0.07392197125256673
$ws.Range("S17").Value = 0.1190965092402464
$ws.Range("F18").Value = 0.01503759398496241
$ws.Range("H18").Value = 0.1353383458646616
$ws.Range("I18").Value = 0.07518796992481203
$ws.Range("J18").Value = 0.5037593984962406
$ws.Range("K18").Value = 0.112781954887218
$ws.Range("M18").Value = 0.007518796992481203
$ws.Range("O18").Value = 0.06015037593984962
$ws.Range("S18").Value = 0.09022556390977443
$ws.Range("F19").Value = 0.01538461538461539
$ws.Range("H19").Value = 0.2085470085470086
$ws.Range("I19").Value = 0.07948717948717948
$ws.Range("J19").Value = 0.3957264957264957
$ws.Range("K19").Value = 0.1136752136752137
$ws.Range("M19").Value = 0.01282051282051282
$ws.Range("O19").Value = 0.07008547008547009
$ws.Range("S19").Value = 0.1042735042735043
